# Auto-generated Excel COM-interop script
# Applies profit-sheet value updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 100
$ws.Range("K8").Value = 300
$ws.Range("M8").Value = -161
# Row 16
$ws.Range("H16").Value = 5500
$ws.Range("J16").Value = 7000
$ws.Range("L16").Value = 7000
$ws.Range("N16").Value = -7460
# Row 80
$ws.Range("H80").Value = 1291.8334
$ws.Range("I80").Value = 1699.5
$ws.Range("J80").Value = 1088
$ws.Range("K80").Value = 5098.5
$ws.Range("L80").Value = 3264
$ws.Range("M80").Value = -4100.5
$ws.Range("N80").Value = -5260
# Row 83
$ws.Range("H83").Value = 1291.8334
$ws.Range("I83").Value = 1699.5
$ws.Range("J83").Value = 1088
$ws.Range("K83").Value = 15295.5
$ws.Range("L83").Value = 9792
$ws.Range("M83").Value = -10303.5
$ws.Range("N83").Value = -19776
# Row 96
$ws.Range("H96").Value = 369.52942
$ws.Range("I96").Value = 274.7857
$ws.Range("K96").Value = 824.3571000000001
$ws.Range("M96").Value = 548.6428999999999
# Row 112
$ws.Range("H112").Value = 2567.7058
$ws.Range("I112").Value = 1652.3334
$ws.Range("J112").Value = 2763.8572
$ws.Range("K112").Value = 4957.0002
$ws.Range("L112").Value = 8291.571599999999
$ws.Range("M112").Value = -3849.0002
$ws.Range("N112").Value = -10507.5716
# Row 113
$ws.Range("H113").Value = 2499.5
$ws.Range("I113").Value = 2499.5
$ws.Range("K113").Value = 2499.5
$ws.Range("M113").Value = 754.5
# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
# Row 132
$ws.Range("H132").Value = 1521.2354
$ws.Range("I132").Value = 1521.2354
$ws.Range("K132").Value = 4563.706200000001
$ws.Range("M132").Value = -2033.706200000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 14922
$ws.Range("I74").Value = 14922
$ws.Range("K74").Value = 14922
$ws.Range("M74").Value = -14048
# Row 77
$ws.Range("H77").Value = 14922
$ws.Range("I77").Value = 14922
$ws.Range("K77").Value = 74610
$ws.Range("M77").Value = -70242
# Row 81
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
# Row 84
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 3487.6
$ws.Range("I64").Value = 594
$ws.Range("J64").Value = 5416.6665
$ws.Range("K64").Value = 594
$ws.Range("L64").Value = 5416.6665
$ws.Range("M64").Value = -369
$ws.Range("N64").Value = -5866.6665
# Row 67
$ws.Range("H67").Value = 3487.6
$ws.Range("I67").Value = 594
$ws.Range("J67").Value = 5416.6665
$ws.Range("K67").Value = 594
$ws.Range("L67").Value = 5416.6665
$ws.Range("M67").Value = 186
$ws.Range("N67").Value = -6976.6665
# Row 94
$ws.Range("H94").Value = 665.625
$ws.Range("I94").Value = 618.7143
$ws.Range("K94").Value = 618.7143
$ws.Range("M94").Value = -167.7143
# Row 102
$ws.Range("H102").Value = 30666.334
$ws.Range("I102").Value = 24999.5
$ws.Range("J102").Value = 42000
$ws.Range("K102").Value = 24999.5
$ws.Range("L102").Value = 42000
$ws.Range("M102").Value = -21754.5
$ws.Range("N102").Value = -48490

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 4192.7144
$ws.Range("I105").Value = 3969.8
$ws.Range("K105").Value = 3969.8
$ws.Range("M105").Value = -2222.8
# Row 132
$ws.Range("H132").Value = 3518.7
$ws.Range("I132").Value = 2479.8
$ws.Range("K132").Value = 7439.400000000001
$ws.Range("M132").Value = -4909.400000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 4510.7144
$ws.Range("J34").Value = 4510.7144
$ws.Range("L34").Value = 13532.1432
$ws.Range("N34").Value = -13700.1432
# Row 60
$ws.Range("H60").Value = 430
$ws.Range("I60").Value = 287.5
$ws.Range("K60").Value = 862.5
$ws.Range("M60").Value = -611.5
# Row 81
$ws.Range("H81").Value = 6123.625
$ws.Range("J81").Value = 5998.5713
$ws.Range("L81").Value = 17995.7139
$ws.Range("N81").Value = -20241.7139
# Row 84
$ws.Range("H84").Value = 6123.625
$ws.Range("J84").Value = 5998.5713
$ws.Range("L84").Value = 53987.14169999999
$ws.Range("N84").Value = -65219.14169999999
# Row 122
$ws.Range("H122").Value = 3477.6956
$ws.Range("I122").Value = 1504
$ws.Range("J122").Value = 3567.4092
$ws.Range("K122").Value = 13536
$ws.Range("L122").Value = 32106.6828
$ws.Range("M122").Value = -11086
$ws.Range("N122").Value = -37006.6828

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 279
$ws.Range("J2").Value = 116
$ws.Range("L2").Value = 116
$ws.Range("N2").Value = -342
# Row 12
$ws.Range("H12").Value = 1500
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = -860
$ws.Range("N12").Value = -2280
# Row 126
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 3000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -7940

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 30000
$ws.Range("J14").Value = 30000
$ws.Range("L14").Value = 30000
$ws.Range("N14").Value = -30344
# Row 22
$ws.Range("H22").Value = 7249.6
$ws.Range("I22").Value = 5348
$ws.Range("J22").Value = 7725
$ws.Range("K22").Value = 5348
$ws.Range("L22").Value = 7725
$ws.Range("M22").Value = -5053
$ws.Range("N22").Value = -8315
# Row 27
$ws.Range("H27").Value = 7249.6
$ws.Range("I27").Value = 5348
$ws.Range("J27").Value = 7725
$ws.Range("K27").Value = 5348
$ws.Range("L27").Value = 7725
$ws.Range("M27").Value = -5241
$ws.Range("N27").Value = -7939
# Row 64
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10450
# Row 67
$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11560
# Row 93
$ws.Range("H93").Value = 1906.8
$ws.Range("I93").Value = 1838.4286
$ws.Range("K93").Value = 1838.4286
$ws.Range("M93").Value = -590.4286
# Row 132
$ws.Range("H132").Value = 4437
$ws.Range("I132").Value = 3767.7144
$ws.Range("K132").Value = 11303.1432
$ws.Range("M132").Value = -8773.143199999999
# Row 136
$ws.Range("H136").Value = 3787.5
$ws.Range("I136").Value = 3787.5
$ws.Range("K136").Value = 11362.5
$ws.Range("M136").Value = -8812.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 88
$ws.Range("H88").Value = 50000
$ws.Range("J88").Value = 50000
$ws.Range("L88").Value = 50000
$ws.Range("N88").Value = -50812
# Row 91
$ws.Range("H91").Value = 50000
$ws.Range("J91").Value = 50000
$ws.Range("L91").Value = 50000
$ws.Range("N91").Value = -52808
# Row 113
$ws.Range("H113").Value = 559.2857
$ws.Range("I113").Value = 569.1667
$ws.Range("K113").Value = 1707.5001
$ws.Range("M113").Value = 462.4999
# Row 126
$ws.Range("H126").Value = 500
$ws.Range("I126").Value = 500
$ws.Range("K126").Value = 1500
$ws.Range("M126").Value = 970
# Row 132
$ws.Range("H132").Value = 2069.3044
$ws.Range("I132").Value = 1841.3158
$ws.Range("K132").Value = 5523.9474
$ws.Range("M132").Value = -2993.9474
